$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row (41) with the latest run log entry, matching the style
# of the existing data rows (copy formatting from row 40).
$newRow = 41

$ws.Range("A40:H40").Copy()
$ws.Range("A$newRow`:H$newRow").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Cells.Item($newRow, 1).Value = "2025-08-21 13:04:19 UTC"
$ws.Cells.Item($newRow, 2).Value = "2025-08-21 18:34:19 IST"
$ws.Cells.Item($newRow, 3).Value = "SKIPPED"
$ws.Cells.Item($newRow, 4).Value = "No change in PDF. Skipping download & Excel update."
$ws.Cells.Item($newRow, 5).Value = "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-21-08-2025.pdf"
$ws.Cells.Item($newRow, 6).Value = ""
$ws.Cells.Item($newRow, 7).Value = 0
$ws.Cells.Item($newRow, 8).Value = ""
